$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column I ("2021") mirrors column H's formatting for rows 4-25.
# Copy each H cell's format into the matching I cell (so the new cells
# reuse the same cell-style index as their H neighbour), then fill in
# the 2021 values where the source diff has them (some rows, the
# "category" rows, stay blank - only the style is copied).

$xlPasteFormats = -4122

$values = @{
    4  = 2021
    5  = 48.5
    7  = 48.8
    8  = 48.2
    10 = 58.2
    11 = 42.4
    12 = 40.7
    14 = 41.5
    15 = 52.6
    17 = 67.1
    18 = 62
    19 = 46.9
    20 = 55.8
    21 = 42.7
    22 = 48.3
    23 = 39.7
    24 = 38.1
    25 = 44.7
}

for ($row = 4; $row -le 25; $row++) {
    $ws.Range("H$row").Copy() | Out-Null
    $ws.Range("I$row").PasteSpecial($xlPasteFormats) | Out-Null

    if ($values.ContainsKey($row)) {
        $ws.Range("I$row").Value = $values[$row]
    }
}

$excel.CutCopyMode = $false

# Collapse the lingering A14:C15 selection back down to the default A1
# cell so the saved sheet view no longer carries the old selection.
$ws.Range("A1").Select() | Out-Null
